$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "45.903.18"
Set-TextValue "E2" "  +3.02%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.449.77"
Set-TextValue "E3" "  +0.77%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "321.29"
Set-TextValue "E5" "  +2.70%  "

# Row 6 - Solana
Set-TextValue "D6" "104.24"
Set-TextValue "E6" "  +2.29%  "

# Row 7 - XRP
Set-TextValue "E7" "  +0.88%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.532"
Set-TextValue "E9" "  +3.80%  "

# Row 10 - Avalanche
Set-TextValue "D10" "35.78"
Set-TextValue "E10" "  +1.36%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  +0.36%  "

# Row 12 - TRON
Set-TextValue "E12" "  -1.82%  "

# Row 13 - Chainlink
Set-TextValue "D13" "18.16"
Set-TextValue "E13" "  -3.32%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.05"
Set-TextValue "E14" "  +1.37%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.832.55"
Set-TextValue "E15" "  +0.81%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.445.94"
Set-TextValue "E16" "  +1.13%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.839"
Set-TextValue "E17" "  -0.08%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "45.773.68"
Set-TextValue "E18" "  +2.87%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "12.51"
Set-TextValue "E19" "  +0.68%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.40"
Set-TextValue "E20" "  -0.03%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0932"
Set-TextValue "E21" "  +2.72%  "

# Row 22 - Litecoin
Set-TextValue "D22" "71.29"
Set-TextValue "E22" "  +3.38%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "246.06"
Set-TextValue "E23" "  +1.98%  "

# Row 24 - ImmutableX
Set-TextValue "D24" "2.35"
Set-TextValue "E24" "  +2.54%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.51"
Set-TextValue "E25" "  +0.86%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "25.93"
Set-TextValue "E26" "  +2.77%  "

# Row 27 - Dai
Set-TextValue "E27" "  -0.01%  "

# Row 28 - was Toncoin, now Cosmos
Set-TextValue "B28" "Cosmos"
Set-TextValue "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "9.69"
Set-TextValue "E28" "  +0.36%  "

# Row 29 - was Cosmos, now Toncoin
Set-TextValue "B29" "Toncoin"
Set-TextValue "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.07"
Set-TextValue "E29" "  -8.73%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "33.68"
Set-TextValue "E30" "  +1.07%  "

# Row 31 - OKB
Set-TextValue "D31" "49.17"
Set-TextValue "E31" "  +1.44%  "

# Row 32 - Kaspa
Set-TextValue "D32" "0.128"
Set-TextValue "E32" "  +4.28%  "

# Row 33 - Celestia
Set-TextValue "D33" "19.87"
Set-TextValue "E33" "  +1.76%  "

# Row 34 - Filecoin
Set-TextValue "E34" "  +2.50%  "

# Row 35 - FirstDigitalUSD
Set-TextValue "E35" "  -0.03%  "

# Row 36 - Hedera
Set-TextValue "E36" "  -0.59%  "

# Row 37 - RenderToken
Set-TextValue "D37" "4.53"
Set-TextValue "E37" "  -0.26%  "

# Row 38 - ARBITRUM
Set-TextValue "E38" "  -0.69%  "

# Row 39 - LidoDAOToken
Set-TextValue "E39" "  +0.72%  "

# Row 40 - Monero
Set-TextValue "D40" "126.34"
Set-TextValue "E40" "  -0.64%  "

# Row 41 - WEMIXToken
Set-TextValue "E41" "  +3.37%  "

# Row 42 - Stellar
Set-TextValue "E42" "  +1.42%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "21.05"
Set-TextValue "E43" "  -4.58%  "

# Row 44 - VeChain
Set-TextValue "E44" "  +1.00%  "

# Row 45 - Maker
Set-TextValue "D45" "1.955.78"
Set-TextValue "E45" "  +0.40%  "

# Row 47 - NEARProtocol
Set-TextValue "E47" "  +0.10%  "

# Row 48 - Stacks
Set-TextValue "D48" "1.84"
Set-TextValue "E48" "  +9.24%  "

# Row 49 - FraxShare
Set-TextValue "D49" "9.08"
Set-TextValue "E49" "  -7.49%  "

# Row 50 - BitcoinSV
Set-TextValue "D50" "77.49"
Set-TextValue "E50" "  +4.89%  "

# Row 51 - THORChain
Set-TextValue "D51" "4.93"
Set-TextValue "E51" "  +5.99%  "
